$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.329.67'
$ws.Range('E2').Value = '  +0.36%  '

$ws.Range('D3').Value = '1.906.85'
$ws.Range('E3').Value = '  +2.61%  '

$ws.Range('E4').Value = '  -0.40%  '

$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.61'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +2.65%  '

$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.668'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  +7.35%  '

$ws.Range('E7').Value = '  -0.37%  '

$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.29'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  -2.27%  '

$ws.Range('E9').Value = '  +6.05%  '

$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.79'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  +12.56%  '

$ws.Range('E11').Value = '  +3.63%  '

$__style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0995'
$ws.Range('D12').Style = $__style
$ws.Range('E12').Value = '  +0.49%  '

$ws.Range('D13').Value = '2.185.31'
$ws.Range('E13').Value = '  +2.70%  '

$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.10'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  +5.19%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.698'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +3.21%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.917.22'
$ws.Range('E16').Value = '  +2.98%  '

$ws.Range('D18').Value = '35.314.29'
$ws.Range('E18').Value = '  +0.34%  '

$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.61'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  +3.89%  '

$ws.Range('E20').Value = '  +4.13%  '

$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '239.54'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  -0.53%  '

$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.49'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  +2.16%  '

$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.84'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  +1.71%  '

$ws.Range('E24').Value = '  -0.45%  '

$ws.Range('E25').Value = '  +1.07%  '

$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.32'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  +22.05%  '

$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.11'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +0.17%  '

$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.45'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  +5.39%  '

$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.42'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  +4.22%  '

$ws.Range('E30').Value = '  +2.62%  '

$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.14'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  +3.47%  '

$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0564'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  +0.58%  '

$ws.Range('B33').Value = 'BinanceUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.02'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  +0.51%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.933'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  +13.79%  '

$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.10'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +2.20%  '

$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.76'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  -3.10%  '

$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.04'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  +0.04%  '

$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.34'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('E39').Value = '  +0.94%  '

$ws.Range('E40').Value = '  +3.18%  '

$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.26'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  +8.66%  '

$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0642'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  +6.53%  '

$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('D44').Value = '1.340.37'
$ws.Range('E44').Value = '  -0.52%  '

$ws.Range('E45').Value = '  +2.58%  '

$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '47.45'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  +36.65%  '

$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.40'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  -0.86%  '

$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.78'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  +1.42%  '

$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.55'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -0.45%  '

$ws.Range('D50').Value = '2.092.39'
$ws.Range('E50').Value = '  +2.29%  '

$ws.Range('E51').Value = '  +3.78%  '
